$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-12 15:58:10.534000 to 2024-03-12 16:56:55.260000"
$ws.Range("B1").Style = "Normal"

$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").Value = 0.04070815972222222

$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = 33.03811083333333

$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = 1691.071816753333

$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = 39.476

$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = 7.414

$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 99

$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 18

$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = 40.78227420008778

$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = 41.46585372989556

$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = 81

$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Custom mode`n74.37%`nEco mode`n20.71%`nSports mode`n1.64%"

$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = 4502.33508

$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = -1738.891328281062

$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = 0.0764918

$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 0.004523068710953785

$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.394

$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 3.094

$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = 0.3000000000000003

$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = 36

$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = 45

$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 9

$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = 61

$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = 60

$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = 60

$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = 50

$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = 95

$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = 0

$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 45

$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 36

$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = 9

$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 53

$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.751019874166667

$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.0000001382980976658347

$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 41

$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = 10.21450141651879

$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 4.321160611438

$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = 2.375393045048411

$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = 5.220883534136546

$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = 10.88073223125058

$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = 11.01771426792441

$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = 51.5457177547399

$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = 4.342953208181564

$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 0

$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0

$ws.Range("B2").NumberFormat = "[hh]:mm:ss"
